$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = @'
Studio: 7
1 bed: 4
2 bed: 3
Percent: 4.0%
'@

$ws.Range("D2").Value = @'
6/16: $2217-2424
($5.1-5.97)
'@

$ws.Range("E2").Value = @'
6/16: $2861-3034
($4.98-5.21)
'@

$ws.Range("F2").Value = @'
6/16: $3510-4300
($4.03-4.21)
'@

$ws.Range("B3").Value = @'
Studio: 4
1 bed: 11
2 bed: 8
Percent: 12.1%
'@

$ws.Range("C3").Value = @'
Save 8 weeks off rent! Minimum 13 month lease term required.
'@

$ws.Range("D3").Value = @'
6/16: $2208-2688
($3.59-4.37)
'@

$ws.Range("E3").Value = @'
6/16: $2566-3496
($2.91-4.16)
'@

$ws.Range("F3").Value = @'
6/16: $3693-7411
($2.68-4.07)
'@

$ws.Range("B4").Value = @'
Studio: 6
1 bed: 18
2 bed: 3
Percent: 6.8%
'@

$ws.Range("D4").Value = @'
6/16: $2050-2225
($3.84-4.85)
'@

$ws.Range("E4").Value = @'
6/16: $2570-3442
($3.04-4.36)
'@

$ws.Range("F4").Value = @'
6/16: $3105-3407
($3.16-3.51)
'@

$ws.Range("B5").Value = @'
Studio: 2
1 bed: 4
2 bed: 6
Percent: 7.5%
'@

$ws.Range("D5").Value = @'
6/16: $2312-2406
($4.61-4.75)
'@

$ws.Range("E5").Value = @'
6/16: $2572-2643
($3.4-3.86)
'@

$ws.Range("F5").Value = @'
6/16: $3441-3508
($3.41-3.51)
'@

$ws.Range("B6").Value = @'
Studio: 6
1 bed: 6
2 bed: 1
Percent: 4.3%
'@

$ws.Range("C6").Value = @'
Up to 1 month off select apartment homes. Offer valid on new leases only. Transfers excluded. [Offer good thru Jun 29, 2024]
'@

$ws.Range("D6").Value = @'
6/16: $1826-2060
($3.55-4.77)
'@

$ws.Range("E6").Value = @'
6/16: $2219-2382
($3.82-4.09)
'@

$ws.Range("F6").Value = @'
6/16: $3386
($3.51)
'@

$ws.Range("B7").Value = @'
Studio: 7
1 bed: 7
2 bed: 6
Percent: 4.5%
'@

$ws.Range("D7").Value = @'
6/16: $1856-2045
($3.69-4.12)
'@

$ws.Range("E7").Value = @'
6/16: $2355-2702
($2.97-3.3)
'@

$ws.Range("F7").Value = @'
6/16: $3119-3275
($2.94-3.08)
'@

$ws.Range("B8").Value = @'
Studio: 7
1 bed: 15
2 bed: 14
Percent: 7.8%
'@

$ws.Range("C8").Value = @'
6 weeks FREE on all apts! Look & lease & receive 2 weeks free *Offers subject to change. Some restrictions may apply. Contact office for details
'@

$ws.Range("D8").Value = @'
6/16: $2695-2915
($3.51-3.8)
'@

$ws.Range("E8").Value = @'
6/16: $2470-3340
($4.01-5.15)
'@

$ws.Range("F8").Value = @'
6/16: $3470-4314
($3.52-4.28)
'@

$ws.Range("B9").Value = @'
Studio: 6
1 bed: 7
Percent: 4.3%
'@

$ws.Range("D9").Value = @'
6/16: $1599-1950
($3.48-4.24)
'@

$ws.Range("E9").Value = @'
6/16: $1825-1999
($2.72-2.98)
'@

$ws.Range("F9").Value = @'
n/a
'@

$ws.Range("B10").Value = @'
Studio: 5
1 bed: 9
2 bed: 4
Percent: 8.2%
'@

$ws.Range("D10").Value = @'
6/16: $1916-2193
($3.09-3.5)
'@

$ws.Range("E10").Value = @'
6/16: $2372-2708
($3.17-3.89)
'@

$ws.Range("F10").Value = @'
6/16: $2808-3441
($2.78-3.18)
'@

$ws.Range("B11").Value = @'
Studio: 3
1 bed: 1
2 bed: 3
Percent: 11.7%
'@

$ws.Range("D11").Value = @'
6/16: $1925-2149
($3.61-4.61)
'@

$ws.Range("E11").Value = @'
6/16: $2199
($3.5)
'@

$ws.Range("F11").Value = @'
6/16: $3021-3200
($3.2-3.66)
'@

$ws.Range("B12").Value = @'
1 bed: 1
2 bed: 5
Percent: 10.7%
'@

$ws.Range("E12").Value = @'
6/16: $2210
($2.93)
'@

$ws.Range("F12").Value = @'
6/16: $2599-2991
($2.93-3.19)
'@

$ws.Range("B13").Value = @'
Studio: 5
1 bed: 6
Percent: 6.8%
'@

$ws.Range("D13").Value = @'
6/16: $1900
($4.22)
'@

$ws.Range("E13").Value = @'
6/16: $2200-2500
($2.93-3.13)
'@

$ws.Range("B14").Value = @'
Studio: 3
1 bed: 1
2 bed: 1
Percent: 4.9%
'@

$ws.Range("D14").Value = @'
6/16: $2245-2295
($3.74-4.01)
'@

$ws.Range("E14").Value = @'
6/16: $2595
($3.33)
'@

$ws.Range("F14").Value = @'
6/16: $3090
($2.81)
'@

$ws.Range("B15").Value = @'
Studio: 1
1 bed: 1
2 bed: 1
Percent: 3.2%
'@

$ws.Range("D15").Value = @'
6/16: $1600
($3.33)
'@

$ws.Range("E15").Value = @'
6/16: $1900
($3.13)
'@

$ws.Range("F15").Value = @'
6/16: $2655
'@

$ws.Range("B16").Value = @'
1 bed: 5
2 bed: 3
Percent: 6.6%
'@

$ws.Range("E16").Value = @'
6/16: $2524-2831
($3.41-3.87)
'@

$ws.Range("F16").Value = @'
6/16: $3370-3563
($3.56-3.83)
'@
